$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1699.7407
$ws.Range("J17").Value = 1339.6154
$ws.Range("L17").Value = 4018.8462
$ws.Range("N17").Value = -4354.8462
$ws.Range("H98").Value = 2639.4211
$ws.Range("I98").Value = 2327.8333
$ws.Range("K98").Value = 2327.8333
$ws.Range("M98").Value = -829.8332999999998
$ws.Range("H122").Value = 2639.4211
$ws.Range("I122").Value = 2327.8333
$ws.Range("K122").Value = 6983.499899999999
$ws.Range("M122").Value = -4533.499899999999
$ws.Range("H129").Value = 1014.2414
$ws.Range("J129").Value = 1009.58
$ws.Range("L129").Value = 3028.74
$ws.Range("N129").Value = -13028.74
$ws.Range("H131").Value = 2345.25
$ws.Range("I131").Value = 631.6667
$ws.Range("K131").Value = 1895.0001
$ws.Range("M131").Value = 3144.9999
$ws.Range("H138").Value = 3384.9333
$ws.Range("I138").Value = 3706.7273
$ws.Range("J138").Value = 2500
$ws.Range("K138").Value = 11120.1819
$ws.Range("L138").Value = 7500
$ws.Range("M138").Value = -5980.1819
$ws.Range("N138").Value = -17780
$ws.Range("H141").Value = 1002267.7
$ws.Range("I141").Value = 1274613.5
$ws.Range("J141").Value = 3666.5
$ws.Range("K141").Value = 3823840.5
$ws.Range("L141").Value = 10999.5
$ws.Range("M141").Value = -3818660.5
$ws.Range("N141").Value = -21359.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1662752.2
$ws.Range("I2").Value = 2327055.5
$ws.Range("J2").Value = 1994.25
$ws.Range("K2").Value = 2327055.5
$ws.Range("L2").Value = 1994.25
$ws.Range("M2").Value = -2326942.5
$ws.Range("N2").Value = -2220.25
$ws.Range("H32").Value = 2165.218
$ws.Range("I32").Value = 1514.7937
$ws.Range("J32").Value = 4897
$ws.Range("K32").Value = 1514.7937
$ws.Range("L32").Value = 4897
$ws.Range("M32").Value = -1227.7937
$ws.Range("N32").Value = -5471
$ws.Range("H43").Value = 38900
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("H61").Value = 3430.4
$ws.Range("I61").Value = 2638
$ws.Range("J61").Value = 4161.846
$ws.Range("K61").Value = 2638
$ws.Range("L61").Value = 4161.846
$ws.Range("M61").Value = -2426
$ws.Range("N61").Value = -4585.846
$ws.Range("H116").Value = 1662752.2
$ws.Range("I116").Value = 2327055.5
$ws.Range("J116").Value = 1994.25
$ws.Range("K116").Value = 2327055.5
$ws.Range("L116").Value = 1994.25
$ws.Range("M116").Value = -2324761.5
$ws.Range("N116").Value = -6582.25
$ws.Range("H132").Value = 2377.7046
$ws.Range("I132").Value = 2374.457
$ws.Range("J132").Value = 2390.3333
$ws.Range("K132").Value = 7123.370999999999
$ws.Range("L132").Value = 7170.999899999999
$ws.Range("M132").Value = -4593.370999999999
$ws.Range("N132").Value = -12230.9999
$ws.Range("H136").Value = 3430.4
$ws.Range("I136").Value = 2638
$ws.Range("J136").Value = 4161.846
$ws.Range("K136").Value = 7914
$ws.Range("L136").Value = 12485.538
$ws.Range("M136").Value = -5364
$ws.Range("N136").Value = -17585.538
$ws.Range("M43").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1662752.2
$ws.Range("I3").Value = 2327055.5
$ws.Range("J3").Value = 1994.25
$ws.Range("K3").Value = 2327055.5
$ws.Range("L3").Value = 1994.25
$ws.Range("M3").Value = -2326941.5
$ws.Range("N3").Value = -2222.25
$ws.Range("H99").Value = 1332.1111
$ws.Range("I99").Value = 1332.1111
$ws.Range("K99").Value = 1332.1111
$ws.Range("M99").Value = 165.8888999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 904.75
$ws.Range("I16").Value = 904.75
$ws.Range("K16").Value = 904.75
$ws.Range("M16").Value = -617.75
$ws.Range("H31").Value = 1533.1305
$ws.Range("J31").Value = 2018.1154
$ws.Range("L31").Value = 2018.1154
$ws.Range("N31").Value = -2608.1154
$ws.Range("H34").Value = 1533.1305
$ws.Range("J34").Value = 2018.1154
$ws.Range("L34").Value = 2018.1154
$ws.Range("N34").Value = -2422.1154
$ws.Range("H113").Value = 904.75
$ws.Range("I113").Value = 904.75
$ws.Range("K113").Value = 904.75
$ws.Range("M113").Value = 1265.25
$ws.Range("H132").Value = 2904.238
$ws.Range("I132").Value = 2714.5
$ws.Range("J132").Value = 3283.7144
$ws.Range("K132").Value = 8143.5
$ws.Range("L132").Value = 9851.143199999999
$ws.Range("M132").Value = -5613.5
$ws.Range("N132").Value = -14911.1432
$ws.Range("H134").Value = 1954.6364
$ws.Range("I134").Value = 943.1429000000001
$ws.Range("K134").Value = 2829.4287
$ws.Range("M134").Value = -294.4287000000004

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 51.133335
$ws.Range("I12").Value = 20.11111
$ws.Range("J12").Value = 97.666664
$ws.Range("K12").Value = 60.33333
$ws.Range("L12").Value = 292.999992
$ws.Range("M12").Value = 112.66667
$ws.Range("N12").Value = -638.999992
$ws.Range("H19").Value = 1500
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("H68").Value = 1849.4897
$ws.Range("J68").Value = 1991.7906
$ws.Range("L68").Value = 5975.3718
$ws.Range("N68").Value = -7597.3718
$ws.Range("H71").Value = 1849.4897
$ws.Range("J71").Value = 1991.7906
$ws.Range("L71").Value = 17926.1154
$ws.Range("N71").Value = -26038.1154
$ws.Range("H98").Value = 456.4
$ws.Range("J98").Value = 452.2857
$ws.Range("L98").Value = 1356.8571
$ws.Range("N98").Value = -4352.8571
$ws.Range("H107").Value = 1344.6666
$ws.Range("J107").Value = 1386.2667
$ws.Range("L107").Value = 4158.800099999999
$ws.Range("N107").Value = -7998.800099999999
$ws.Range("H131").Value = 12839684
$ws.Range("J131").Value = 20179.676
$ws.Range("L131").Value = 60539.028
$ws.Range("N131").Value = -70619.02799999999
$ws.Range("N19").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2780955
$ws.Range("J126").Value = 2936.889
$ws.Range("L126").Value = 8810.667000000001
$ws.Range("N126").Value = -13750.667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3912
$ws.Range("I7").Value = 2718.5
$ws.Range("K7").Value = 2718.5
$ws.Range("M7").Value = -2606.5
$ws.Range("H16").Value = 4907
$ws.Range("I16").Value = 5090.8
$ws.Range("K16").Value = 5090.8
$ws.Range("M16").Value = -4920.8
$ws.Range("H126").Value = 3912
$ws.Range("I126").Value = 2718.5
$ws.Range("K126").Value = 8155.5
$ws.Range("M126").Value = -5685.5
$ws.Range("H132").Value = 3697.2942
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 3865.875
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 11597.625
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -16657.625
$ws.Range("H136").Value = 5038.8823
$ws.Range("I136").Value = 2997.2
$ws.Range("K136").Value = 8991.599999999999
$ws.Range("M136").Value = -6441.599999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 45397.516
$ws.Range("I122").Value = 48652.926
$ws.Range("K122").Value = 145958.778
$ws.Range("M122").Value = -143508.778
$ws.Range("H126").Value = 4757.4287
$ws.Range("I126").Value = 3861.0667
$ws.Range("K126").Value = 11583.2001
$ws.Range("M126").Value = -9113.2001
$ws.Range("H136").Value = 23150928
$ws.Range("I136").Value = 32682108
$ws.Range("J136").Value = 3778.5715
$ws.Range("K136").Value = 98046324
$ws.Range("L136").Value = 11335.7145
$ws.Range("M136").Value = -98043774
$ws.Range("N136").Value = -16435.7145
